# This workbook ("Pais" sheet) is a COVID-19 stats scrape: column A holds the
# country name, B..H hold the daily metrics (Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes), and the
# whole table (rows 4..219) is kept sorted by column B ("Casos totales")
# descending. This run's source update (i) refreshes the timestamp in A1,
# (ii) refreshes USA/Brazil case numbers, and (iii) refreshes the numbers for
# Guinea, Gabon, Santa Lucia, Groenlandia and San Bartolome - which nudges
# those countries (and the neighbours they swap past) into new rank
# positions in the sorted table. We therefore just overwrite, cell by cell,
# every row whose displayed country/number actually differs afterwards.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = 'Datos actualizados a 24 de Mayo de 2020 a las 01:35'

# Row 4
$ws.Cells.Item(4,2).Value = 1666736
$ws.Cells.Item(4,3).Value = 21757
$ws.Cells.Item(4,4).Value = 446866
$ws.Cells.Item(4,5).Value = 1121197
$ws.Cells.Item(4,7).Value = 1026
$ws.Cells.Item(4,8).Value = 98673

# Row 5
$ws.Cells.Item(5,4).Value = 142587
$ws.Cells.Item(5,5).Value = 182798

# Row 75
$ws.Cells.Item(75,1).Value = 'Guinea'
$ws.Cells.Item(75,2).Value = 3176
$ws.Cells.Item(75,3).Value = 109
$ws.Cells.Item(75,4).Value = 1631
$ws.Cells.Item(75,5).Value = 1525
$ws.Cells.Item(75,7).Value = 1
$ws.Cells.Item(75,8).Value = 20

# Row 76
$ws.Cells.Item(76,1).Value = 'Uzbekistan'
$ws.Cells.Item(76,2).Value = 3115
$ws.Cells.Item(76,3).Value = 87
$ws.Cells.Item(76,4).Value = 2532
$ws.Cells.Item(76,5).Value = 570
$ws.Cells.Item(76,8).Value = 13

# Row 89
$ws.Cells.Item(89,1).Value = 'Gabon'
$ws.Cells.Item(89,2).Value = 1934
$ws.Cells.Item(89,3).Value = 206
$ws.Cells.Item(89,4).Value = 459
$ws.Cells.Item(89,5).Value = 1463
$ws.Cells.Item(89,8).Value = 12

# Row 90
$ws.Cells.Item(90,1).Value = 'Cuba'
$ws.Cells.Item(90,2).Value = 1931
$ws.Cells.Item(90,3).Value = 15
$ws.Cells.Item(90,4).Value = 1671
$ws.Cells.Item(90,5).Value = 179
$ws.Cells.Item(90,8).Value = 81

# Row 91
$ws.Cells.Item(91,1).Value = 'Estonia'
$ws.Cells.Item(91,2).Value = 1821
$ws.Cells.Item(91,3).Value = 14
$ws.Cells.Item(91,4).Value = 1526
$ws.Cells.Item(91,5).Value = 231
$ws.Cells.Item(91,8).Value = 64

# Row 92
$ws.Cells.Item(92,1).Value = 'El Salvador'
$ws.Cells.Item(92,2).Value = 1819
$ws.Cells.Item(92,3).Value = 94
$ws.Cells.Item(92,4).Value = 570
$ws.Cells.Item(92,5).Value = 1216
$ws.Cells.Item(92,8).Value = 33

# Row 93
$ws.Cells.Item(93,1).Value = 'Islandia'
$ws.Cells.Item(93,2).Value = 1804
$ws.Cells.Item(93,3).Value = 1
$ws.Cells.Item(93,4).Value = 1791
$ws.Cells.Item(93,5).Value = 3
$ws.Cells.Item(93,8).Value = 10

# Row 198
$ws.Cells.Item(198,1).Value = 'Santa Lucia'
$ws.Cells.Item(198,4).Value = 18
$ws.Cells.Item(198,8).Value = 0

# Row 199
$ws.Cells.Item(199,1).Value = 'Belice'
$ws.Cells.Item(199,4).Value = 16
$ws.Cells.Item(199,8).Value = 2

# Row 200
$ws.Cells.Item(200,1).Value = 'Nueva Caledonia'

# Row 210
$ws.Cells.Item(210,1).Value = 'Groenlandia'
$ws.Cells.Item(210,4).Value = 11
$ws.Cells.Item(210,8).Value = 0

# Row 211
$ws.Cells.Item(211,1).Value = 'Montserrat'
$ws.Cells.Item(211,4).Value = 10
$ws.Cells.Item(211,8).Value = 1

# Row 215
$ws.Cells.Item(215,1).Value = 'San Bartolome'

# Row 216
$ws.Cells.Item(216,1).Value = 'Bonaire, San Eustaquio y Saba'
